$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "hang san xuat & loai hang" - the ngayChungTu/ngayThanhToan/tuNgay/denNgay
# sample dates in row 2 are switched from numeric date serials to literal
# formatted-text values (so the template shows the raw placeholder text
# instead of a real date), and the now-wider text needs a bit less column
# width than the old numeric display.
$ws.Range("B2").Value = "03-08-2018T00:00:00"
$ws.Range("F2").Value = "03-08-2018T00:00:00"
$ws.Range("I2").Value = "03-08-2018T00:00:00"
$ws.Range("J2").Value = "03-11-2018T00:00:00"

$ws.Columns("B").ColumnWidth = 19.140625
$ws.Columns("F").ColumnWidth = 19.140625
$ws.Range("I:J").ColumnWidth = 19.140625
